$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")
$ws.Activate()

# Copy style from G7 (header style) to H7:I7
$ws.Range("G7").Copy()
$ws.Range("H7:I7").PasteSpecial(-4122)

# Copy style from B8 (data style) to H8:I10
$ws.Range("B8").Copy()
$ws.Range("H8:I10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Set values
$ws.Range("H7").Value = "AlarmLoadingDetail"
$ws.Range("I7").Value = "StandbyLoadingDetail"
$ws.Range("H8").Value = "Battery Alarm (A)"
$ws.Range("I8").Value = "Battery Standby (A)"
$ws.Range("H9").Value = "Battery Alarm (A)"
$ws.Range("I9").Value = "Battery Standby (A)"
$ws.Range("H10").Value = "Battery Alarm (A)"
$ws.Range("I10").Value = "Battery Standby (A)"

# Update selection to H10:I10 with active cell H10
$ws.Range("H10:I10").Select()

Write-Output "done"
